$d = $word.ActiveDocument

# Helper: find a unique paragraph by its (unique) text substring and return the Paragraph object
function Get-ParaByText($text) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $text"
        return $null
    }
    return $r.Paragraphs(1)
}

# Helper: replace the whole text of a range found by search with new text (merges runs, drops proofErr)
function Replace-Text($oldText, $newText) {
    $ok = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Host "REPLACE FAILED FOR: $oldText"
    }
    return $ok
}

# 1) Heading1: merge "How to integrate this " + "Instagram" + " Page Designer component"
Replace-Text "How to integrate this Instagram Page Designer component" "How to integrate this Instagram Page Designer component"

# 2) "Import " + "MetaData" + " (in zip)" -> single run, proofErr removed
Replace-Text "Import MetaData (in zip)" "Import MetaData (in zip)"

# 3) "Add the following information to the YouTube settings in the Site Preferences" -> "Import Jobs (in zip)"
Replace-Text "Add the following information to the YouTube settings in the Site Preferences" "Import Jobs (in zip)"

# 4) Delete paragraph "apiKey - copied from previous steps"
$p = Get-ParaByText("apiKey - copied from previous steps")
if ($p -ne $null) {
    $p.Range.Delete()
}

# 5) Insert new paragraph "Import Custom Object (in zip)" before "Add files from cartridge..."
$p = Get-ParaByText("Add files from cartridge into appropriate paths on the site as laid out in the sample cartridge.")
$newRange = $p.Range.InsertParagraphBefore()
$p2 = Get-ParaByText("Add files from cartridge into appropriate paths on the site as laid out in the sample cartridge.")
$prevPara = $p2.Previous()
$prevPara.Range.Text = "Import Custom Object (in zip)"

# 6) Remove the empty paragraph right after "Add files from cartridge..."
$p = Get-ParaByText("Add files from cartridge into appropriate paths on the site as laid out in the sample cartridge.")
$next = $p.Next()
if ($next.Range.Text.Trim().Length -eq 0) {
    $next.Range.Delete()
}

# 7) "De" + "A" + "uth" + " Callback and Deletion URLs are the same. " -> single run, proofErr removed
# (the proofErr spellStart sits as the very first child of this paragraph, so anchor the Find
#  to the end of the previous paragraph via the ^p paragraph-mark code to ensure the replaced
#  range fully encloses it and the markers get dropped)
Replace-Text "OAuth Redirect URL the same.^pDeAuth Callback and Deletion URLs are the same. " "OAuth Redirect URL the same.^pDeAuth Callback and Deletion URLs are the same. "

# 8) "In " + "Instagram" (Heading2) -> "In Instagram"
Replace-Text "In Instagram" "In Instagram"

# 9) Merge the access_token sentence into a single run, dropping proofErr wraps
Replace-Text "The access_token received from the second authentication should be placed in the Business Manager Settings for the appropriate site.  (LongToken)" "The access_token received from the second authentication should be placed in the Business Manager Settings for the appropriate site.  (LongToken)"

# 10) Remove the empty paragraph right after "Run the job, and the Instagram Custom Object..."
$p = Get-ParaByText("Run the job, and the Instagram Custom Object will be populated with the most recent posts")
$next = $p.Next()
if ($next.Range.Text.Trim().Length -eq 0) {
    $next.Range.Delete()
}

# 11) Merge "Drag the "Instagram" component..." into a single run
Replace-Text "Drag the “Instagram” component to the area of the page you want to place the feed." "Drag the “Instagram” component to the area of the page you want to place the feed."

# 12) Merge "Save and preview. Your feed should show!" into a single run
Replace-Text "Save and preview. Your feed should show!" "Save and preview. Your feed should show!"

# 13) Replace "If you want to add text around the feed..." text with "There is a title component that is available for use. "
Replace-Text "If you want to add text around the feed, you can customize this component, or use the Text component above or below for more control. " "There is a title component that is available for use. "

# 14) Delete the (now orphaned) empty paragraph that had ind left=720 hanging=360
$p = Get-ParaByText("There is a title component that is available for use.")
$next = $p.Next()
if ($next.Range.Text.Trim().Length -eq 0) {
    $next.Range.Delete()
}

Write-Host "All edits applied."
